# Apply translation-sheet update described in the commit:
#   "Bug with translation of '&'"
#
# Net effect on xl/worksheets/sheet1.xml:
#   1. A brand new row is inserted at row 75 ("Initial and Final Surveillance
#      Diagnosis" / TBT / new), pushing every following row down by one.
#   2. The row that (after the shift) lands at 145 - which used to hold
#      "Susceptible & Intermediate are always combined in this visualisation
#      of co-resistances." - is overwritten with a corrected translation
#      using "and" instead of "&", flagged as "new".
#   3. A final new row (187) is appended, re-adding the old "&" phrased
#      string but flagged as "deleted" (it is being retired in favour of the
#      "and" phrasing added above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 75 (everything from old row 75 onward
#    shifts down by one row).
$ws.Rows.Item(75).Insert()

$ws.Range("A75").Value = "Initial and Final Surveillance Diagnosis"
$ws.Range("B75").Value = "TBT"
$ws.Range("C75").Value = "new"

# 2. Fix the "Susceptible & Intermediate ..." translation that is now
#    (after the shift above) sitting at row 145, replacing the "&" with
#    "and" and flagging the row as new.
$ws.Range("A145").Value = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
$ws.Range("B145").Value = "TBT"
$ws.Range("C145").Value = "new"

# 3. Append a row re-recording the old "&" wording, flagged as deleted.
$ws.Range("A187").Value = "Susceptible & Intermediate are always combined in this visualisation of co-resistances."
$ws.Range("B187").Value = "TBT"
$ws.Range("C187").Value = "deleted"
